$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the shared formulas in row 2 (C2:O2) with plain literal values.
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 300
$ws.Range("G2").Value = 400
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 111
$ws.Range("N2").Value = 50
$ws.Range("O2").Value = 5

# Restore the selection to match the reverted state.
$ws.Range("F10").Select()
